# "added summer winter filter"
# Adds a new "season" column (F) to the meals sheet, classifying each
# recipe as "both", "winter" or "summer".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: season -------------------------------------------------
$ws.Range("F1").Value = "season"

$ws.Range("F2").Value  = "both"    # baked salmon and veggies
$ws.Range("F3").Value  = "both"    # burgers
$ws.Range("F4").Value  = "both"    # chicken cheese wraps
$ws.Range("F5").Value  = "both"    # chicken mango chutney rice
$ws.Range("F6").Value  = "both"    # chicken_wings
$ws.Range("F7").Value  = "winter"  # couscous with stew
$ws.Range("F8").Value  = "both"    # curry
$ws.Range("F9").Value  = "both"    # fish with fries
$ws.Range("F10").Value = "winter"  # fiskeboller
$ws.Range("F11").Value = "winter"  # fiskegrateng
$ws.Range("F12").Value = "winter"  # kjottboller
$ws.Range("F13").Value = "both"    # miso soup
$ws.Range("F14").Value = "both"    # nachos
$ws.Range("F15").Value = "summer"  # paella
$ws.Range("F16").Value = "both"    # pasta
$ws.Range("F17").Value = "both"    # pie
$ws.Range("F18").Value = "winter"  # stamppot
$ws.Range("F19").Value = "both"    # stir fry noodles
$ws.Range("F20").Value = "both"    # gnocchi

# --- Widen column F so the season labels are readable ----------------------
$ws.Columns("F").ColumnWidth = 33.6

# --- Rows whose text now wraps to a second line need a taller row height --
$ws.Rows(11).RowHeight = 28.8
$ws.Rows(15).RowHeight = 28.8
$ws.Rows(18).RowHeight = 28.8

# --- Restore the view: scroll the frozen pane down a bit and leave the
#     selection on the new filter cell --------------------------------------
$ws.Range("F21").Select()
